$d = $word.ActiveDocument

# 1) Bullet: "Event sourced .Net 7 application..." -> "Contributed to the design, development..."
$d.Content.Find.Execute(
    "Event sourced .Net 7 application with realtime processing of Kafka messages, a REST API and back office tool.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Contributed to the design, development, and maintenance of event-sourced, gamification-oriented applications featuring real-time Kafka message processing, RESTful APIs, and Blazor-based back-office tools",
    2) | Out-Null

# 2) Bullet: "Performance optimizations leading to..." -> "Optimized Kafka message processing..."
$d.Content.Find.Execute(
    "Performance optimizations leading to near 1000% increase in the number of Kafka messages the application can process per second",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Optimized Kafka message processing, resulting in a nearly 1000% increase in throughput for one of the core applications.",
    2) | Out-Null

# 3) Bullet: "Conducting and evaluating interviews..." -> "Developed Proofs of Concept..."
$d.Content.Find.Execute(
    "Conducting and evaluating interviews for potential new recruits",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Developed Proofs of Concept for new systems and features to evaluate feasibility and drive innovation.",
    2) | Out-Null

# 4) Bullet: "Proofs of Concept for new systems and/or functionality" -> "Collaborated with Enterprise Architect..."
$d.Content.Find.Execute(
    "Proofs of Concept for new systems and/or functionality",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Collaborated with Enterprise Architect and Engineering Manager to ensure architectural and regulatory compliance.",
    2) | Out-Null

# 5) Insert a brand-new bullet right after the one above, carrying the (reworded)
#    interview-related sentence, re-using the same list/paragraph formatting.
$rng = $d.Content
$rng.Find.Execute(
    "Collaborated with Enterprise Architect and Engineering Manager to ensure architectural and regulatory compliance.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$origPara = $rng.Paragraphs(1)
$rng.InsertParagraphAfter()
$newPara = $origPara.Next()
$newPara.Range.Text = "Conducted and assessed technical interviews to support recruitment and team growth."

# 6) Skills line update
$d.Content.Find.Execute(
    "C#, .Net 5-9, Asp.Net Core, Blazor, Docker, GitHub Actions, Kafka, SQL Server, PostgreSQL",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "C#, .Net 5-9, .Net Aspire, Asp.Net Core, Microsoft Orleans, Blazor, Docker, GitHub Actions, Kafka, SQL Server, PostgreSQL",
    2) | Out-Null
